# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" worksheet (fund-holding detail, same layout as the
# existing 2021-Q3 / 2021-Q4 sheets) positioned right after "2021-Q4" and
# right before "总计", and prepends a matching "2022-Q1" summary row to the
# "总计" sheet.

$wb = $excel.ActiveWorkbook

# Reference sheet used purely as a formatting donor (header row / column-A
# "index" style) so the new sheet's look matches its siblings exactly.
$donor = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# 1. New "2022-Q1" worksheet, inserted between "2021-Q4" and "总计".
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Add($null, $donor)
$q1.Name = "2022-Q1"

# Pull over the header-row / row-index cell styles from the "2021-Q4" sheet.
$donor.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

$donor.Range("A2:A4").Copy()
$q1.Range("A2:A4").PasteSpecial(-4122)

# Header row.
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

function Set-TextValue($cell, $text) {
    # Force text storage (no scientific/float coercion, leading zeros kept)
    # for numeric-looking strings, then drop the quote-prefix style Excel
    # applies on entry so the cell ends up style-less, matching the source.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2 — 970023 天风天盈一年定期开放混合
$q1.Cells.Item(2, 1).Value = 0
Set-TextValue $q1.Cells.Item(2, 2) "970023"
$q1.Cells.Item(2, 3).Value = "天风天盈一年定期开放混合"
Set-TextValue $q1.Cells.Item(2, 4) "2.97"
Set-TextValue $q1.Cells.Item(2, 5) "39.47"
Set-TextValue $q1.Cells.Item(2, 6) "4.92"
Set-TextValue $q1.Cells.Item(2, 7) "0.1461"
$q1.Cells.Item(2, 8).Value = 3

# Row 3 — 008033 中加科盈混合A
$q1.Cells.Item(3, 1).Value = 1
Set-TextValue $q1.Cells.Item(3, 2) "008033"
$q1.Cells.Item(3, 3).Value = "中加科盈混合A"
Set-TextValue $q1.Cells.Item(3, 4) "7.58"
Set-TextValue $q1.Cells.Item(3, 5) "23.16"
Set-TextValue $q1.Cells.Item(3, 6) "1.20"
Set-TextValue $q1.Cells.Item(3, 7) "0.0910"
$q1.Cells.Item(3, 8).Value = 7

# Row 4 — 008034 中加科盈混合C
$q1.Cells.Item(4, 1).Value = 2
Set-TextValue $q1.Cells.Item(4, 2) "008034"
$q1.Cells.Item(4, 3).Value = "中加科盈混合C"
Set-TextValue $q1.Cells.Item(4, 4) "1.16"
Set-TextValue $q1.Cells.Item(4, 5) "23.16"
Set-TextValue $q1.Cells.Item(4, 6) "1.20"
Set-TextValue $q1.Cells.Item(4, 7) "0.0139"
$q1.Cells.Item(4, 8).Value = 7

# ---------------------------------------------------------------------
# 2. "总计" sheet: prepend a 2022-Q1 summary row above the existing ones.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Cells.Clear()

$donor.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)

$donor.Range("A2:A4").Copy()
$total.Range("A2:A4").PasteSpecial(-4122)

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 3
$total.Cells.Item(2, 4).Value = 0.25

$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(3, 2).Value = "2021-Q4"
$total.Cells.Item(3, 3).Value = 9
$total.Cells.Item(3, 4).Value = 0.61

$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(4, 2).Value = "2021-Q3"
$total.Cells.Item(4, 3).Value = 6
$total.Cells.Item(4, 4).Value = 0.67
